$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.831958762886598
$ws.Range("B3").Value = 0.8330869608491775
$ws.Range("B4").Value = 0.831958762886598
$ws.Range("B5").Value = 0.8364940290651653
$ws.Range("B6").Value = 0.7039220174718617
